$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 400
$ws.Range("I7").Value = 400
$ws.Range("K7").Value = 400
$ws.Range("M7").Value = -288

$ws.Range("H14").Value = 400
$ws.Range("I14").Value = 400
$ws.Range("K14").Value = 400
$ws.Range("M14").Value = -209

$ws.Range("H107").Value = 1400.2174
$ws.Range("I107").Value = 1709.3334
$ws.Range("K107").Value = 1709.3334
$ws.Range("M107").Value = 210.6666

$ws.Range("H132").Value = 8508.9
$ws.Range("I132").Value = 9223.75
$ws.Range("J132").Value = 5649.5
$ws.Range("K132").Value = 27671.25
$ws.Range("L132").Value = 16948.5
$ws.Range("M132").Value = -25141.25
$ws.Range("N132").Value = -22008.5

$ws.Range("H141").Value = 4123.4546
$ws.Range("I141").Value = 2193.923
$ws.Range("J141").Value = 6910.5557
$ws.Range("K141").Value = 6581.768999999999
$ws.Range("L141").Value = 20731.6671
$ws.Range("M141").Value = -1401.768999999999
$ws.Range("N141").Value = -31091.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2018775
$ws.Range("I32").Value = 2533324.5
$ws.Range("J32").Value = 21112.412
$ws.Range("K32").Value = 2533324.5
$ws.Range("L32").Value = 21112.412
$ws.Range("M32").Value = -2533037.5
$ws.Range("N32").Value = -21686.412

$ws.Range("H74").Value = 10640823
$ws.Range("I74").Value = 1489.3214
$ws.Range("J74").Value = 26319842
$ws.Range("K74").Value = 1489.3214
$ws.Range("L74").Value = 26319842
$ws.Range("M74").Value = -615.3214
$ws.Range("N74").Value = -26321590

$ws.Range("H77").Value = 10640823
$ws.Range("I77").Value = 1489.3214
$ws.Range("J77").Value = 26319842
$ws.Range("K77").Value = 7446.607
$ws.Range("L77").Value = 131599210
$ws.Range("M77").Value = -3078.607
$ws.Range("N77").Value = -131607946

$ws.Range("H112").Value = 43092.875
$ws.Range("J112").Value = 43092.875
$ws.Range("L112").Value = 43092.875
$ws.Range("N112").Value = -46046.875

$ws.Range("H132").Value = 1926843.5
$ws.Range("I132").Value = 2825.6667
$ws.Range("K132").Value = 8477.000100000001
$ws.Range("M132").Value = -5947.000100000001

$ws.Range("H138").Value = 77504.336
$ws.Range("J138").Value = 77504.336
$ws.Range("L138").Value = 77504.336
$ws.Range("N138").Value = -87784.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2481.5898
$ws.Range("I134").Value = 2350.5454
$ws.Range("J134").Value = 3202.3333
$ws.Range("K134").Value = 7051.6362
$ws.Range("L134").Value = 9606.999899999999
$ws.Range("M134").Value = -4516.6362
$ws.Range("N134").Value = -14676.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 17183.455
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 17183.455
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 17183.455
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -17407.455

$ws.Range("H31").Value = 4685.757
$ws.Range("I31").Value = 1115.6
$ws.Range("J31").Value = 7889.7437
$ws.Range("K31").Value = 1115.6
$ws.Range("L31").Value = 7889.7437
$ws.Range("M31").Value = -820.5999999999999
$ws.Range("N31").Value = -8479.743699999999

$ws.Range("H34").Value = 4685.757
$ws.Range("I34").Value = 1115.6
$ws.Range("J34").Value = 7889.7437
$ws.Range("K34").Value = 1115.6
$ws.Range("L34").Value = 7889.7437
$ws.Range("M34").Value = -913.5999999999999
$ws.Range("N34").Value = -8293.743699999999

$ws.Range("H132").Value = 6412608.5
$ws.Range("I132").Value = 2161.4666
$ws.Range("J132").Value = 15154127
$ws.Range("K132").Value = 6484.399800000001
$ws.Range("L132").Value = 45462381
$ws.Range("M132").Value = -3954.399800000001
$ws.Range("N132").Value = -45467441

$ws.Range("H134").Value = 10006520
$ws.Range("I134").Value = 15632812
$ws.Range("J134").Value = 4222.1113
$ws.Range("K134").Value = 46898436
$ws.Range("L134").Value = 12666.3339
$ws.Range("M134").Value = -46895901
$ws.Range("N134").Value = -17736.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4119574.5
$ws.Range("I4").Value = 15556136
$ws.Range("J4").Value = 2412.04
$ws.Range("K4").Value = 46668408
$ws.Range("L4").Value = 7236.12
$ws.Range("M4").Value = -46668296
$ws.Range("N4").Value = -7460.12

$ws.Range("H122").Value = 2254.7048
$ws.Range("I122").Value = 349.65384
$ws.Range("J122").Value = 3669.8857
$ws.Range("K122").Value = 3146.88456
$ws.Range("L122").Value = 33028.9713
$ws.Range("M122").Value = -696.88456
$ws.Range("N122").Value = -37928.9713

$ws.Range("H125").Value = 2888.5715
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2888.5715
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 8665.7145
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -18505.7145

$ws.Range("H131").Value = 3591.2327
$ws.Range("J131").Value = 4804.968
$ws.Range("L131").Value = 14414.904
$ws.Range("N131").Value = -24494.904

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224

$ws.Range("H28").Value = 19999
$ws.Range("J28").Value = 19999
$ws.Range("L28").Value = 19999
$ws.Range("N28").Value = -20383

$ws.Range("H70").Value = 15356.333
$ws.Range("I70").Value = 15356.333
$ws.Range("K70").Value = 15356.333
$ws.Range("M70").Value = -15086.333

$ws.Range("H73").Value = 15356.333
$ws.Range("I73").Value = 15356.333
$ws.Range("K73").Value = 15356.333
$ws.Range("M73").Value = -14420.333

$ws.Range("H132").Value = 2957.5386
$ws.Range("I132").Value = 2068.2222
$ws.Range("J132").Value = 4958.5
$ws.Range("K132").Value = 6204.6666
$ws.Range("L132").Value = 14875.5
$ws.Range("M132").Value = -3674.6666
$ws.Range("N132").Value = -19935.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8168
$ws.Range("J2").Value = 8168
$ws.Range("L2").Value = 8168
$ws.Range("N2").Value = -8392

$ws.Range("H7").Value = 5596.773
$ws.Range("I7").Value = 5427.8423
$ws.Range("J7").Value = 6666.6665
$ws.Range("K7").Value = 5427.8423
$ws.Range("L7").Value = 6666.6665
$ws.Range("M7").Value = -5315.8423
$ws.Range("N7").Value = -6890.6665

$ws.Range("H126").Value = 5596.773
$ws.Range("I126").Value = 5427.8423
$ws.Range("J126").Value = 6666.6665
$ws.Range("K126").Value = 16283.5269
$ws.Range("L126").Value = 19999.9995
$ws.Range("M126").Value = -13813.5269
$ws.Range("N126").Value = -24939.9995

$ws.Range("H132").Value = 3175.2856
$ws.Range("I132").Value = 2714.6843
$ws.Range("J132").Value = 3722.25
$ws.Range("K132").Value = 8144.0529
$ws.Range("L132").Value = 11166.75
$ws.Range("M132").Value = -5614.0529
$ws.Range("N132").Value = -16226.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 502500000
$ws.Range("J2").Value = 1000000000
$ws.Range("L2").Value = 1000000000
$ws.Range("N2").Value = -1000000224

$ws.Range("H28").Value = 42579.715
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 42579.715
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 42579.715
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -43275.715

$ws.Range("H126").Value = 841.35
$ws.Range("I126").Value = 841.35
$ws.Range("K126").Value = 2524.05
$ws.Range("M126").Value = -54.05000000000018

$ws.Range("H132").Value = 9412243
$ws.Range("I132").Value = 4824.4165
$ws.Range("J132").Value = 15353771
$ws.Range("K132").Value = 14473.2495
$ws.Range("L132").Value = 46061313
$ws.Range("M132").Value = -11943.2495
$ws.Range("N132").Value = -46066373

$ws.Range("H136").Value = 3050.7754
$ws.Range("I136").Value = 2770.742
$ws.Range("J136").Value = 3533.0557
$ws.Range("K136").Value = 8312.226000000001
$ws.Range("L136").Value = 10599.1671
$ws.Range("M136").Value = -5762.226000000001
$ws.Range("N136").Value = -15699.1671
